# Emissies per km2 excl puntbrongegevens - 2000 NOx
# Add the (previously removed) 2000 data marker row back under the header,
# and leave the selection where Excel would land after entering it (A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2000

$ws.Range("A3").Select() | Out-Null
